$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 30 (pushes old row 30 -> 31, totals 31 -> 32, footer 32 -> 33)
$ws.Rows(30).Insert()

# Fill in the new row 30 values (new item: "معجون سيجنال 25 مل")
$ws.Range("A30").Value = 24
$ws.Range("C30").Value = "معجون سيجنال 25 مل"
$ws.Range("H30").Value = "18:0"
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = "20.00"
$ws.Range("P30").Value = "20.0000"
$ws.Range("Q30").Value = 16

# Update totals (now on row 32) to add the new item's sale value
$ws.Range("P32").Value = 852.09

# Update the timestamp in the footer row (now row 33)
$ws.Range("A33").Value = "Sunday, 14 September, 2025 1:23 PM"
